# This script rewrites the PL_specific_ion_cfg lookup table on Sheet1 to match
# the updated configuration: adds Canopy IDE support data (new PC HG / PE HG
# "whitelist" NEG rows for PC and PE), fixes a ppm-related duplicate-index bug
# in the PA/PC/PE blocks, and removes an erroneous duplicate PS "serine" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the used range completely (contents + stray styled-but-empty cells);
# we will rewrite the whole table since rows are reordered/shifted/added.
$ws.UsedRange.Clear()

$rows = 43
$cols = 8
$arr = New-Object 'object[,]' $rows,$cols

# Row 1: CLASS / LABEL
$arr[0,0] = "CLASS"
$arr[0,1] = "TYPE"
$arr[0,2] = "EXACTMASS"
$arr[0,3] = "FORMULA"
$arr[0,4] = "CHARGE_MODE"
$arr[0,5] = "PR_CHARGE"
$arr[0,6] = "LABEL"
$arr[0,7] = "REMARKS"

# Row 2: PA / PA:-98
$arr[1,0] = "PA"
$arr[1,1] = "NL"
$arr[1,2] = 97.976898000000006
$arr[1,3] = "H3O4P"
$arr[1,4] = "NEG"
$arr[1,5] = "[M-H]-"
$arr[1,6] = "PA:-98"
$arr[1,7] = "PA HG"

# Row 3: PA / PA:-98
$arr[2,0] = "PA"
$arr[2,1] = "NL"
$arr[2,2] = 97.976898000000006
$arr[2,3] = "H3O4P"
$arr[2,4] = "POS"
$arr[2,5] = "[M+H]+"
$arr[2,6] = "PA:-98"
$arr[2,7] = "PA HG"

# Row 4: PA / PA:-98
$arr[3,0] = "PA"
$arr[3,1] = "NL"
$arr[3,2] = 97.976898000000006
$arr[3,3] = "H3O4P"
$arr[3,4] = "POS"
$arr[3,5] = "[M+NH4]+"
$arr[3,6] = "PA:-98"
$arr[3,7] = "PA HG"

# Row 6: PC / PC:168
$arr[5,0] = "PC"
$arr[5,1] = "FRAG"
$arr[5,2] = 168.04580000000001
$arr[5,3] = "C4H11O4NP-"
$arr[5,4] = "NEG"
$arr[5,5] = "[M+HCOO]-"
$arr[5,6] = "PC:168"
$arr[5,7] = "deprotonated demethylated PC"

# Row 7: PC / PC:224
$arr[6,0] = "PC"
$arr[6,1] = "FRAG"
$arr[6,2] = 224.06878699999999
$arr[6,3] = "C7H15O5NP-"
$arr[6,4] = "NEG"
$arr[6,5] = "[M+HCOO]-"
$arr[6,6] = "PC:224"
$arr[6,7] = "demethylated PC dehydrated glycerol ester"

# Row 8: PC / PC:242
$arr[7,0] = "PC"
$arr[7,1] = "FRAG"
$arr[7,2] = 242.079352
$arr[7,3] = "C7H17O6NP-"
$arr[7,4] = "NEG"
$arr[7,5] = "[M+HCOO]-"
$arr[7,6] = "PC:242"
$arr[7,7] = "demethylated PC glycerol ester"

# Row 9: PC / PC:-60
$arr[8,0] = "PC"
$arr[8,1] = "NL"
$arr[8,2] = 60.021129999999999
$arr[8,3] = "C2H4O2"
$arr[8,4] = "NEG"
$arr[8,5] = "[M+HCOO]-"
$arr[8,6] = "PC:-60"
$arr[8,7] = "methyl formate"

# Row 10: PC / PC:-183
$arr[9,0] = "PC"
$arr[9,1] = "NL"
$arr[9,2] = 183.066047
$arr[9,3] = "C5H14NO4P"
$arr[9,4] = "NEG"
$arr[9,5] = "[M+HCOO]-"
$arr[9,6] = "PC:-183"
$arr[9,7] = "PC HG"

# Row 11: PC / PC:-59
$arr[10,0] = "PC"
$arr[10,1] = "NL"
$arr[10,2] = 59.073498999999998
$arr[10,3] = "C3H9N"
$arr[10,4] = "POS"
$arr[10,5] = "[M+H]+"
$arr[10,6] = "PC:-59"
$arr[10,7] = "(CH3)3N"

# Row 12: PC / PC:-59
$arr[11,0] = "PC"
$arr[11,1] = "NL"
$arr[11,2] = 59.073498999999998
$arr[11,3] = "C3H9N"
$arr[11,4] = "POS"
$arr[11,5] = "[M+NH4]+"
$arr[11,6] = "PC:-59"
$arr[11,7] = "(CH3)3N"

# Row 13: PC / PC:183
$arr[12,0] = "PC"
$arr[12,1] = "NL"
$arr[12,2] = 183.066047
$arr[12,3] = "C5H14NO4P"
$arr[12,4] = "POS"
$arr[12,5] = "[M+H]+"
$arr[12,6] = "PC:183"
$arr[12,7] = "PC HG"

# Row 14: PC / PC:183
$arr[13,0] = "PC"
$arr[13,1] = "NL"
$arr[13,2] = 183.066047
$arr[13,3] = "C5H14NO4P"
$arr[13,4] = "POS"
$arr[13,5] = "[M+NH4]+"
$arr[13,6] = "PC:183"
$arr[13,7] = "PC HG"

# Row 15: PC / PC:184
$arr[14,0] = "PC"
$arr[14,1] = "FRAG"
$arr[14,2] = 184.07387199999999
$arr[14,3] = "C5H15NO4P+"
$arr[14,4] = "POS"
$arr[14,5] = "[M+H]+"
$arr[14,6] = "PC:184"
$arr[14,7] = "PC HG [M+H]+"

# Row 16: PC / PC:184
$arr[15,0] = "PC"
$arr[15,1] = "FRAG"
$arr[15,2] = 184.07387199999999
$arr[15,3] = "C5H15NO4P+"
$arr[15,4] = "POS"
$arr[15,5] = "[M+NH4]+"
$arr[15,6] = "PC:184"
$arr[15,7] = "PC HG [M+H]+"

# Row 18: PE / PE:140
$arr[17,0] = "PE"
$arr[17,1] = "FRAG"
$arr[17,2] = 140.01127199999999
$arr[17,3] = "C2H7O4NP-"
$arr[17,4] = "NEG"
$arr[17,5] = "[M-H]-"
$arr[17,6] = "PE:140"
$arr[17,7] = "deprotonated phosphoethanolamine"

# Row 19: PE / PE:196
$arr[18,0] = "PE"
$arr[18,1] = "FRAG"
$arr[18,2] = 196.037487
$arr[18,3] = "C5H11O5NP-"
$arr[18,4] = "NEG"
$arr[18,5] = "[M-H]-"
$arr[18,6] = "PE:196"
$arr[18,7] = "deprotonated doubly dehydrated glycerol phosphocholine (dilyso-h2o)"

# Row 20: PE / PE:-141
$arr[19,0] = "PE"
$arr[19,1] = "NL"
$arr[19,2] = 141.01909699999999
$arr[19,3] = "C2H8NO4P"
$arr[19,4] = "NEG"
$arr[19,5] = "[M-H]-"
$arr[19,6] = "PE:-141"
$arr[19,7] = "PE HG"

# Row 21: PE / PE:-43
$arr[20,0] = "PE"
$arr[20,1] = "NL"
$arr[20,2] = 43.042198999999997
$arr[20,3] = "C2H5N"
$arr[20,4] = "NEG"
$arr[20,5] = "[M-H]-"
$arr[20,6] = "PE:-43"
$arr[20,7] = "PE HG part"

# Row 22: PE / PE:142
$arr[21,0] = "PE"
$arr[21,1] = "FRAG"
$arr[21,2] = 142.02692200000001
$arr[21,3] = "C2H9NO4P+"
$arr[21,4] = "POS"
$arr[21,5] = "[M+H]+"
$arr[21,6] = "PE:142"
$arr[21,7] = "PE HG [M+H]+"

# Row 23: PE / PE:142
$arr[22,0] = "PE"
$arr[22,1] = "FRAG"
$arr[22,2] = 142.02692200000001
$arr[22,3] = "C2H9NO4P+"
$arr[22,4] = "POS"
$arr[22,5] = "[M+NH4]+"
$arr[22,6] = "PE:142"
$arr[22,7] = "PE HG [M+H]+"

# Row 24: PE / PE:-43
$arr[23,0] = "PE"
$arr[23,1] = "NL"
$arr[23,2] = 43.042198999999997
$arr[23,3] = "C2H5N"
$arr[23,4] = "POS"
$arr[23,5] = "[M+H]+"
$arr[23,6] = "PE:-43"
$arr[23,7] = "PE HG part"

# Row 25: PE / PE:-43
$arr[24,0] = "PE"
$arr[24,1] = "NL"
$arr[24,2] = 43.042198999999997
$arr[24,3] = "C2H5N"
$arr[24,4] = "POS"
$arr[24,5] = "[M+NH4]+"
$arr[24,6] = "PE:-43"
$arr[24,7] = "PE HG part"

# Row 27: PG / PG:171
$arr[26,0] = "PG"
$arr[26,1] = "FRAG"
$arr[26,2] = 171.005853
$arr[26,3] = "C3H8O6P-"
$arr[26,4] = "NEG"
$arr[26,5] = "[M-H]-"
$arr[26,6] = "PG:171"
$arr[26,7] = "phosphoglycerol"

# Row 28: PG / PG:153
$arr[27,0] = "PG"
$arr[27,1] = "FRAG"
$arr[27,2] = 152.99528799999999
$arr[27,3] = "C3H6O5P-"
$arr[27,4] = "NEG"
$arr[27,5] = "[M-H]-"
$arr[27,6] = "PG:153"
$arr[27,7] = "phosphoglycerol - water"

# Row 29: PG / PG:-172
$arr[28,0] = "PG"
$arr[28,1] = "NL"
$arr[28,2] = 172.013678
$arr[28,3] = "C3H9O6P"
$arr[28,4] = "NEG"
$arr[28,5] = "[M-H]-"
$arr[28,6] = "PG:-172"
$arr[28,7] = "PG HG"

# Row 30: PG / PG:-172
$arr[29,0] = "PG"
$arr[29,1] = "NL"
$arr[29,2] = 172.013678
$arr[29,3] = "C3H9O6P"
$arr[29,4] = "POS"
$arr[29,5] = "[M+H]+"
$arr[29,6] = "PG:-172"
$arr[29,7] = "PG HG"

# Row 31: PG / PG:-172
$arr[30,0] = "PG"
$arr[30,1] = "NL"
$arr[30,2] = 172.013678
$arr[30,3] = "C3H9O6P"
$arr[30,4] = "POS"
$arr[30,5] = "[M+NH4]+"
$arr[30,6] = "PG:-172"
$arr[30,7] = "PG HG"

# Row 33: PI / PI:241
$arr[32,0] = "PI"
$arr[32,1] = "FRAG"
$arr[32,2] = 241.01133300000001
$arr[32,3] = "C6H10O8P -"
$arr[32,4] = "NEG"
$arr[32,5] = "[M-H]-"
$arr[32,6] = "PI:241"
$arr[32,7] = "phosphoinositol"

# Row 34: PI / PI:-162
$arr[33,0] = "PI"
$arr[33,1] = "NL"
$arr[33,2] = 162.05282399999999
$arr[33,3] = "C6H10O5"
$arr[33,4] = "NEG"
$arr[33,5] = "[M-H]-"
$arr[33,6] = "PI:-162"
$arr[33,7] = "inositol"

# Row 35: PI / PI:-162
$arr[34,0] = "PI"
$arr[34,1] = "NL"
$arr[34,2] = 162.05282399999999
$arr[34,3] = "C6H10O5"
$arr[34,4] = "POS"
$arr[34,5] = "[M+H]+"
$arr[34,6] = "PI:-162"
$arr[34,7] = "inositol"

# Row 36: PI / PI:-162
$arr[35,0] = "PI"
$arr[35,1] = "NL"
$arr[35,2] = 162.05282399999999
$arr[35,3] = "C6H10O5"
$arr[35,4] = "POS"
$arr[35,5] = "[M+NH4]+"
$arr[35,6] = "PI:-162"
$arr[35,7] = "inositol"

# Row 38: PS / PS:184
$arr[37,0] = "PS"
$arr[37,1] = "FRAG"
$arr[37,2] = 184.001102
$arr[37,3] = "C3H7NO6P-"
$arr[37,4] = "NEG"
$arr[37,5] = "[M-H]-"
$arr[37,6] = "PS:184"
$arr[37,7] = "phosphoserine"

# Row 39: PS / PS:-87
$arr[38,0] = "PS"
$arr[38,1] = "NL"
$arr[38,2] = 87.032028999999994
$arr[38,3] = "C3H5NO2"
$arr[38,4] = "NEG"
$arr[38,5] = "[M-H]-"
$arr[38,6] = "PS:-87"
$arr[38,7] = "serine"

# Row 40: PS / PS:-87
$arr[39,0] = "PS"
$arr[39,1] = "NL"
$arr[39,2] = 87.032028999999994
$arr[39,3] = "C3H5NO2"
$arr[39,4] = "POS"
$arr[39,5] = "[M+H]+"
$arr[39,6] = "PS:-87"
$arr[39,7] = "serine"

# Row 41: PS / PS:-87
$arr[40,0] = "PS"
$arr[40,1] = "NL"
$arr[40,2] = 87.032028999999994
$arr[40,3] = "C3H5NO2"
$arr[40,4] = "POS"
$arr[40,5] = "[M+NH4]+"
$arr[40,6] = "PS:-87"
$arr[40,7] = "serine"

# Row 42: PS / PS:186
$arr[41,0] = "PS"
$arr[41,1] = "FRAG"
$arr[41,2] = 186.016752
$arr[41,3] = "C3H9NO6P+"
$arr[41,4] = "POS"
$arr[41,5] = "[M+H]+"
$arr[41,6] = "PS:186"
$arr[41,7] = "phosphoserine"

# Row 43: PS / PS:186
$arr[42,0] = "PS"
$arr[42,1] = "FRAG"
$arr[42,2] = 186.016752
$arr[42,3] = "C3H9NO6P+"
$arr[42,4] = "POS"
$arr[42,5] = "[M+NH4]+"
$arr[42,6] = "PS:186"
$arr[42,7] = "phosphoserine"

$ws.Range("A1:H43").Value2 = $arr

# Restore the active cell/selection recorded for this sheet
$ws.Range("H14").Select()
